# Fruta / hortaliza, semanal
# Insert a new weekly record as row 448, shifting existing rows 448:560 down to 449:561.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 448 (pushes old 448..560 down to 449..561)
$ws.Rows.Item(448).Insert()

# Populate the new row 448 with the new data point
$ws.Range("A448").Value = 5
$ws.Range("B448").Value = "Macroferia Regional de Talca"
$ws.Range("C448").Value = "Maule"
$ws.Range("D448").Value = 44932
$ws.Range("E448").Value = 7
$ws.Range("F448").Value = 100112043
$ws.Range("G448").Value = "Pepino ensalada"
$ws.Range("H448").Value = "Sin especificar"
$ws.Range("I448").Value = "Primera"
$ws.Range("J448").Value = 500
$ws.Range("K448").Value = 10000
$ws.Range("L448").Value = 10000
$ws.Range("M448").Value = 10000
$ws.Range("N448").Value = "$/caja 80 unidades"
$ws.Range("O448").Value = "Región del Maule"
$ws.Range("P448").Value = 125
$ws.Range("Q448").Value = 80
$ws.Range("R448").Value = "Hortaliza"
